# Apply the updates described by the diff: the four data rows (2-5) were
# re-shuffled (row2<->row4 and row3<->row5 for the Id/Ost/Nord/Publik
# kommentar fields), the Ost/Nord coordinates were rounded to whole
# numbers, and the Starttid/Sluttid (Z/AB) columns were cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 111936770
$ws.Range("Q2").Value = 489837
$ws.Range("R2").Value = 7087463
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").Value = "ringhack färska"

# --- Row 3 ---
$ws.Range("A3").Value = 111936772
$ws.Range("Q3").Value = 489837
$ws.Range("R3").Value = 7087471
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "ringhack färska"

# --- Row 4 ---
$ws.Range("A4").Value = 111936769
$ws.Range("Q4").Value = 489838
$ws.Range("R4").Value = 7087500
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").Value = "ringhack gamla"

# --- Row 5 ---
$ws.Range("A5").Value = 111936773
$ws.Range("Q5").Value = 490003
$ws.Range("R5").Value = 7087487
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").Value = "ringhack gamla"
